$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text NumberFormat for Price cells whose new value would otherwise
# be auto-parsed as a Number by Excel (losing exact text, e.g. "0.560" -> 0.56).
# The sheet stores all Price/Volume cells as text, so these must stay text too.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "56.482.97"
$ws.Range("E2").Value = "  +10.30%  "
$ws.Range("D3").Value = "3.254.96"
$ws.Range("E3").Value = "  +6.19%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "397.99"
$ws.Range("E5").Value = "  +2.61%  "
$ws.Range("D6").Value = "111.12"
$ws.Range("E6").Value = "  +8.80%  "
$ws.Range("D7").Value = "0.560"
$ws.Range("E7").Value = "  +4.57%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.621"
$ws.Range("E9").Value = "  +7.38%  "
$ws.Range("D10").Value = "39.39"
$ws.Range("E10").Value = "  +7.18%  "
$ws.Range("D11").Value = "0.0959"
$ws.Range("E11").Value = "  +13.01%  "
$ws.Range("E12").Value = "  +2.23%  "
$ws.Range("D13").Value = "3.752.09"
$ws.Range("E13").Value = "  +5.85%  "
$ws.Range("D14").Value = "19.27"
$ws.Range("E14").Value = "  +5.49%  "
$ws.Range("E15").Value = "  +5.90%  "
$ws.Range("D16").Value = "3.256.74"
$ws.Range("E16").Value = "  +6.07%  "
$ws.Range("E17").Value = "  +5.44%  "
$ws.Range("D18").Value = "11.02"
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("D19").Value = "56.449.88"
$ws.Range("E19").Value = "  +10.21%  "
$ws.Range("D20").Value = "3.32"
$ws.Range("E20").Value = "  +4.33%  "
$ws.Range("E21").Value = "  +9.40%  "
$ws.Range("D22").Value = "13.06"
$ws.Range("E22").Value = "  +6.51%  "
$ws.Range("D23").Value = "299.09"
$ws.Range("E23").Value = "  +13.15%  "
$ws.Range("D24").Value = "75.33"
$ws.Range("E24").Value = "  +8.11%  "
$ws.Range("E25").Value = "  +3.85%  "
$ws.Range("D26").Value = "8.15"
$ws.Range("E26").Value = "  +3.21%  "
$ws.Range("D27").Value = "28.26"
$ws.Range("E27").Value = "  +5.18%  "
$ws.Range("E28").Value = "  +3.91%  "
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("E30").Value = "  +4.38%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").Value = "0.111"
$ws.Range("E32").Value = "  +6.38%  "
$ws.Range("D33").Value = "11.09"
$ws.Range("E33").Value = "  +6.56%  "
$ws.Range("D34").Value = "36.87"
$ws.Range("E34").Value = "  +3.59%  "
$ws.Range("D35").Value = "0.0488"
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("D37").Value = "51.67"
$ws.Range("E37").Value = "  +3.36%  "
$ws.Range("D38").Value = "3.12"
$ws.Range("E38").Value = "  +25.74%  "
$ws.Range("D39").Value = "3.53"
$ws.Range("E39").Value = "  +5.39%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "17.63"
$ws.Range("E41").Value = "  +6.84%  "
$ws.Range("D42").Value = "134.32"
$ws.Range("E42").Value = "  +2.93%  "
$ws.Range("D43").Value = "1.94"
$ws.Range("E43").Value = "  +6.26%  "
$ws.Range("E44").Value = "  +4.77%  "
$ws.Range("D45").Value = "3.98"
$ws.Range("E45").Value = "  +7.05%  "
$ws.Range("D46").Value = "0.286"
$ws.Range("E46").Value = "  -3.42%  "
$ws.Range("D47").Value = "22.24"
$ws.Range("E47").Value = "  +2.52%  "
$ws.Range("E48").Value = "  +56.54%  "
$ws.Range("D49").Value = "2.150.52"
$ws.Range("E49").Value = "  +4.05%  "
$ws.Range("D50").Value = "2.08"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").Value = "2.42"
$ws.Range("E51").Value = "  -4.00%  "
